$d = $word.ActiveDocument

function Assert-ParaText($para, $expected, $label) {
    $actual = $para.Range.Text -replace "[\x0d\x07]+$", ""
    if ($actual -ne $expected) {
        throw "Assertion failed for '$label': expected [$expected] but found [$actual]"
    }
}

# ------------------------------------------------------------------
# 1) Colour the "Đáp án: ĐĐSĐ" run that immediately precedes "Câu 1:"
#    (there are two identical "Đáp án: ĐĐSĐ" paragraphs in the doc –
#    this is the second one, paragraph #63, 1-based) and delete the
#    now-redundant empty paragraph that used to follow it.
# ------------------------------------------------------------------
$pAnswer1 = $d.Paragraphs.Item(63)
Assert-ParaText $pAnswer1 "Đáp án: ĐĐSĐ" "answer-before-Cau1"
$rngColor = $pAnswer1.Range
$rngColor.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark
$rngColor.Font.Color = 13382400       # 0x0033CC (BGR-encoded: 00|33|CC)

$pEmpty = $d.Paragraphs.Item(64)
Assert-ParaText $pEmpty "" "empty-before-Cau1"
$pEmpty.Range.Delete()

# ------------------------------------------------------------------
# 2) Strip the trailing " (đ)"/" (s)" markers from the four
#    true/false statements under "Câu 1:" and add the missing period
#    for the "Câu 2:" statement c).
# ------------------------------------------------------------------
$d.Content.Find.Execute("a) Trường idNhacsi phải được đặt kiểu dữ liệu là INT. (đ)", $false, $false, $false, $false, $false, $true, 1, $false, "a) Trường idNhacsi phải được đặt kiểu dữ liệu là INT.", 2) | Out-Null
$d.Content.Find.Execute("b) Trường tenNhacsi nên có kiểu dữ liệu là CHAR(255). (s)", $false, $false, $false, $false, $false, $true, 1, $false, "b) Trường tenNhacsi nên có kiểu dữ liệu là CHAR(255).", 2) | Out-Null
$d.Content.Find.Execute("c) Chọn AUTO_INCREMENT cho trường idNhacsi để tự động tăng giá trị. (đ)", $false, $false, $false, $false, $false, $true, 1, $false, "c) Chọn AUTO_INCREMENT cho trường idNhacsi để tự động tăng giá trị.", 2) | Out-Null
$d.Content.Find.Execute("d) Bỏ chọn ALLOW NULL cho idNhacsi để đảm bảo không có giá trị rỗng. (đ)", $false, $false, $false, $false, $false, $true, 1, $false, "d) Bỏ chọn ALLOW NULL cho idNhacsi để đảm bảo không có giá trị rỗng.", 2) | Out-Null

# ------------------------------------------------------------------
# 3) Insert a new paragraph "Đáp án: ĐSĐĐ" right after item d) of
#    "Câu 1:" (and before "Câu 2:").
# ------------------------------------------------------------------
$pD1 = $d.Paragraphs.Item(68)
Assert-ParaText $pD1 "d) Bỏ chọn ALLOW NULL cho idNhacsi để đảm bảo không có giá trị rỗng." "item-d-Cau1"
$pD1.Range.InsertParagraphAfter()
$pNew1 = $d.Paragraphs.Item(69)
$rngNew1 = $pNew1.Range
$rngNew1.MoveEnd(1, -1) | Out-Null
$rngNew1.Text = "Đáp án: ĐSĐĐ"

# ------------------------------------------------------------------
# 4) Strip the trailing " (đ)"/" (s)" markers from the four
#    true/false statements under "Câu 2:" (item c) additionally gets
#    a proper closing period).
# ------------------------------------------------------------------
$d.Content.Find.Execute("a) Sử dụng phím Ctrl+Insert để thêm trường mới. (đ)", $false, $false, $false, $false, $false, $true, 1, $false, "a) Sử dụng phím Ctrl+Insert để thêm trường mới.", 2) | Out-Null
$d.Content.Find.Execute("b) Tên trường mới có tên mặc định là Column1. (s)", $false, $false, $false, $false, $false, $true, 1, $false, "b) Tên trường mới có tên mặc định là Column1.", 2) | Out-Null
$d.Content.Find.Execute("c) Với HeidiSQL, tên trường không phân biệt chữ hoa và chữ thường (đ)", $false, $false, $false, $false, $false, $true, 1, $false, "c) Với HeidiSQL, tên trường không phân biệt chữ hoa và chữ thường.", 2) | Out-Null
$d.Content.Find.Execute("d) Tên trường mới có thể đặt trùng tên với các trường đã có. (s)", $false, $false, $false, $false, $false, $true, 1, $false, "d) Tên trường mới có thể đặt trùng tên với các trường đã có.", 2) | Out-Null

# ------------------------------------------------------------------
# 5) The document's final paragraph is already an empty paragraph
#    (right before the sectPr) — give it the run text "Đáp án: ĐSĐS"
#    instead of inserting a brand-new paragraph.
# ------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($lastIndex)
Assert-ParaText $pLast "" "final-empty-paragraph"
$rngLast = $pLast.Range
$rngLast.MoveEnd(1, -1) | Out-Null
$rngLast.Text = "Đáp án: ĐSĐS"
$rngLast2 = $pLast.Range
$rngLast2.MoveEnd(1, -1) | Out-Null
$rngLast2.Font.Name = "Archivo Narrow"
$rngLast2.Font.NameBi = "Archivo Narrow"
